$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): numeric 1..9 -> text labels "P1".."P9" ---
$ws.Range("B1").Value = "P1"
$ws.Range("C1").Value = "P2"
$ws.Range("D1").Value = "P3"
$ws.Range("E1").Value = "P4"
$ws.Range("F1").Value = "P5"
$ws.Range("G1").Value = "P6"
$ws.Range("H1").Value = "P7"
$ws.Range("I1").Value = "P8"
$ws.Range("J1").Value = "P9"

# --- Data row (row 2): "-" -> "NA" in column E ---
$ws.Range("E2").Value = "NA"

# --- Header row formatting: centered, bold-ish font with explicit black color ---
$ws.Range("A1:J1").ClearFormats()
$ws.Range("A1:J1").Font.Color = 0
$ws.Range("A1:J1").HorizontalAlignment = -4108
$ws.Range("A1:J1").VerticalAlignment = -4108

# --- Data row: drop wrap text, restore auto row height ---
$ws.Range("A2:J2").WrapText = $false
$ws.Rows(2).AutoFit()

# --- Column widths: best-fit based on content ---
$ws.Columns("A:J").AutoFit()

# --- Selection moves from L2 to F2 ---
$ws.Range("F2").Select() | Out-Null
